$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("E4").Value = "1. bugfix/nsdsd" + $nl + "2. feature/remsv1.24" + $nl + "3. feature/v.4.5" + $nl + "4. feature/v1.33" + $nl + "5. hotfix/vv222" + $nl + "6. main"
$ws.Range("G4").Value = "1. 4.5"
$ws.Range("H4").Value = "1. 4.5"
